$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "JianHua Tissue Co., Ltd."
$ws.Range("D3").Value = "Quanzhou Blossom Trading Co., Ltd."
$ws.Range("D4").Value = "Sichuan Petrochemical Yashi Paper Co., Ltd."
$ws.Range("D5").Value = "Shenzhen Telling Commodity Co., Ltd."
$ws.Range("D6").Value = "Qingdao Wellpaper Industrial Co., Ltd."
$ws.Range("D7").Value = "Jiaxing Jdl Paper Products Co., Ltd."
$ws.Range("D8").Value = "Guangxi Mashan Shengsheng Paper Co., Ltd."
$ws.Range("D9").Value = "Dalian Weimei House Ware Co., Ltd."
$ws.Range("D10").Value = "Guangdong Union Eco-Packaging Co., Ltd."
$ws.Range("D11").Value = "Joylife Industry (Dongguan) Co., Ltd."
$ws.Range("A12").Value = "Di lusso in rilievo tovaglioli di carta biodegradabile di bambù cena del tessuto Logo personalizzato stampato per i rivenditori di marca"
$ws.Range("B12").Value = "10 €"
$ws.Range("C12").Value = "Ordine minimo: 1.800 parti"
$ws.Range("D12").Value = "Henrich (shandong) Health Technology Co., Ltd."
$ws.Range("D13").Value = "Shenzhen Telling Commodity Co., Ltd."
$ws.Range("D14").Value = "Huzhou Anji Jiahui Import And Export Co., Ltd."
$ws.Range("D15").Value = "Foshan Bao Shi Jie Hygiene Supplies Co., Ltd."
$ws.Range("D16").Value = "Huzhou Yaojin Nonwoven Technology Co., Ltd."
$ws.Range("D17").Value = "Bright Paper Co., Ltd."
$ws.Range("D18").Value = "Hebei Mountain Environmental Protection Technology Co., Ltd."
$ws.Range("D19").Value = "Dalian Huayufei International Trade Co., Ltd."
$ws.Range("D20").Value = "Qingdao Dongfang Jiarui Int'l Co., Ltd."
$ws.Range("D21").Value = "Hangzhou Mingxuan Sanitary Products Co., Ltd."
$ws.Range("D22").Value = "Hefei Green Way Tableware Co., Ltd."
$ws.Range("D23").Value = "Xiamen Qiaodou Daily Commodity Co., Ltd."
$ws.Range("D24").Value = "Baoding Suiqian Trading Co., Ltd."
$ws.Range("D25").Value = "Quanzhou Blossom Trading Co., Ltd."
$ws.Range("D26").Value = "JianHua Tissue Co., Ltd."
$ws.Range("D27").Value = "Shenzhen Telling Commodity Co., Ltd."
$ws.Range("A28").Value = "Tovaglioli di carta monouso tessuti molli di bambù biodegradabili Cocktail 2strati tavolo festa di nozze"
$ws.Range("B28").Value = "0,0087 €"
$ws.Range("C28").Value = "Ordine minimo: 50 parti"
$ws.Range("D28").Value = "Henrich (shandong) Health Technology Co., Ltd."
$ws.Range("B29").Value = "1,04 €"
$ws.Range("D29").Value = "Qingdao Wellpaper Industrial Co., Ltd."
$ws.Range("D31").Value = "Huzhou Anji Jiahui Import And Export Co., Ltd."
$ws.Range("D32").Value = "Foshan Bao Shi Jie Hygiene Supplies Co., Ltd."
$ws.Range("D33").Value = "JianHua Tissue Co., Ltd."
$ws.Range("A34").Value = "Commercio all'ingrosso personalizzato biodegradabile monouso compostabile in legno di bambù stoviglie posate forchetta e cucchiaio coltello set di tovaglioli"
$ws.Range("B34").Value = "0,0173-0,0432 €"
$ws.Range("C34").Value = "Ordine minimo: 5.000 insiemi"
$ws.Range("D34").Value = "Shenzhen Telling Commodity Co., Ltd."
$ws.Range("E34").Value = "5.0"
$ws.Range("A35").Value = "All'ingrosso 100% di bambù tovaglioli di carta marrone usa e getta Eco Friendly tovaglioli di tessuto per la festa di nozze eventi di cena"
$ws.Range("B35").Value = "0,0777-0,1209 €"
$ws.Range("C35").Value = "Ordine minimo: 3.000 fogli"
$ws.Range("D35").Value = "Quanzhou Blossom Trading Co., Ltd."
$ws.Range("E35").Value = "4.5"
$ws.Range("A36").Value = "Set di posate in bambù biodegradabile eco-friendly forchette usa Pre-arrotolate tovaglioli per viaggi di nozze cibo carta verde"
$ws.Range("B36").Value = "0,0087-0,0259 €"
$ws.Range("C36").Value = "Ordine minimo: 100 parti"
$ws.Range("D36").Value = "Huzhou Anji Jiahui Import And Export Co., Ltd."
$ws.Range("E36").Value = "4.8"
$ws.Range("A37").Value = "Eco-friendly 3Ply monouso stoviglie compostabili in fibra di bambù tovaglioli di tessuto per eventi tovaglioli di carta"
$ws.Range("B37").Value = "10 €"
$ws.Range("C37").Value = "Ordine minimo: 50 parti"
$ws.Range("D37").Value = "Henrich (shandong) Health Technology Co., Ltd."
$ws.Range("E37").Value = "4.9"
$ws.Range("A38").Value = "Tovagliolo monouso in polpa di bambù 2 strati di tessuto Eco Friendly biodegradabile"
$ws.Range("B38").Value = "0,5612-1,13 €"
$ws.Range("C38").Value = "Ordine minimo: 2 sacchi"
$ws.Range("D38").Value = "Quanzhou Blossom Trading Co., Ltd."
$ws.Range("E38").Value = "4.5"
$ws.Range("A39").Value = "Set di posate in bambù biodegradabile Pre-arrotolato tovaglioli usa e getta forchette coltelli cucchiai di carta stoviglie per viaggi di nozze"
$ws.Range("B39").Value = "0,0087-0,0259 €"
$ws.Range("C39").Value = "Ordine minimo: 100 parti"
$ws.Range("D39").Value = "Huzhou Anji Jiahui Import And Export Co., Ltd."
$ws.Range("E39").Value = "4.8"
$ws.Range("A40").Value = "Produttore personalizzato all'ingrosso eco-friendly 8 volte 2Ply polpa di bambù marrone chiaro ristorante goffratura Hotel bar tovaglioli"
$ws.Range("B40").Value = "0,1554-0,1986 €"
$ws.Range("C40").Value = "Ordine minimo: 800.000 parti"
$ws.Range("D40").Value = "Hangzhou Mingxuan Sanitary Products Co., Ltd."
$ws.Range("E40").Value = "5.0"
$ws.Range("A41").Value = "Ingrosso stile classico usa e getta eco-friendly pasta di bambù e bagassa insalatiera scatola da pranzo biodegradabile scatola di carta scatole di imballaggio"
$ws.Range("B41").Value = "0,0518 €"
$ws.Range("C41").Value = "Ordine minimo: 50.000 parti"
$ws.Range("D41").Value = "Yiwu Shuangtong Daily Necessities Co., Ltd. Branch"
$ws.Range("E41").ClearContents()
$ws.Range("A42").Value = "Vendita diretta della fabbrica eco-friendly polpa di bambù bianco tovaglioli da Cocktail per cena tovaglioli di carta"
$ws.Range("B42").Value = "0,1554-0,1986 €"
$ws.Range("C42").Value = "Ordine minimo: 1.000 pacchetti"
$ws.Range("D42").Value = "Joylife Industry (Dongguan) Co., Ltd."
$ws.Range("E42").Value = "5.0"
$ws.Range("A43").Value = "Durable and Eco-Friendly 2Ply Virgin Wood Pulp Napkins Colorful Bamboo Sheets in Bag for Dinner Use"
$ws.Range("B43").Value = "0,0087-0,0259 €"
$ws.Range("C43").Value = "Ordine minimo: 50.000 parti"
$ws.Range("D43").Value = "Foshan Bao Shi Jie Hygiene Supplies Co., Ltd."
$ws.Range("A44").Value = "Bacchette di legno di bambù usa e getta all'ingrosso biodegradabili carta kraft ambientale o stoviglie confezionate indipendenti OPP"
$ws.Range("C44").Value = "Ordine minimo: 10.000 parti"
$ws.Range("D44").Value = "Suqian Green Wooden Products Co., Ltd."
$ws.Range("E44").Value = "4.4"
$ws.Range("A45").Value = "Set di posate di bambù 50 confezioni confezionate singolarmente forchette compostabili coltelli tovaglioli campeggio matrimoni-6.7 `"usa e getta"
$ws.Range("B45").Value = "0,0259-0,0777 €"
$ws.Range("C45").Value = "Ordine minimo: 50 parti"
$ws.Range("D45").Value = "Guangzhou Gorlando Commodity Co., Ltd."
$ws.Range("E45").Value = "5.0"
$ws.Range("A46").Value = "Vassoio in fibra riciclata diretta in fabbrica, imballaggio in pasta di bambù biodegradabile, scatola di carta, imballaggio compostabile"
$ws.Range("B46").Value = "0,0691-0,0864 €"
$ws.Range("C46").Value = "Ordine minimo: 1.000 parti"
$ws.Range("D46").Value = "Dongguan Kinyi Packaging Technology Co., Ltd."
$ws.Range("A47").Value = "Tovagliolo Pre arrotolato ecologico e posate di bambù Set di posate compostabili posate confezionate biodegradabili per feste"
$ws.Range("B47").Value = "10 €"
$ws.Range("C47").Value = "Ordine minimo: 5.000 insiemi"
$ws.Range("D47").Value = "Hunan Gianty New Material Technology Co., Ltd."
$ws.Range("A48").ClearContents()
$ws.Range("B48").Value = "0,0864 €"
$ws.Range("C48").Value = "Ordine minimo: 10.000 insiemi"
$ws.Range("D48").Value = "Dalian Huayufei International Trade Co., Ltd."
$ws.Range("E48").Value = "4.8"

# Remove the now-obsolete last row (its data has been absorbed into row 48)
$ws.Rows.Item(49).Delete()
